# Refresh the ligand/receptor TPM-derived metrics (Vcan-Itgb1) with the
# updated TPM values from the new run of the scripts.
# Columns E,F,K,L (cell counts / detection rates) are unaffected; only the
# average/total expression values (G,H,M,N) and everything derived from
# them (I,J,O,P,Q,R,S,T specificity + weight columns) change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.37596266666667
$ws.Range("H2").Value = 34.127888
$ws.Range("I2").Value = 0.05604480707695051
$ws.Range("J2").Value = 0.05604480707695052
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 1814.279904248206
$ws.Range("R2").Value = 16328.51913823386
$ws.Range("S2").Value = 0.0167197637038891
$ws.Range("T2").Value = 0.0167197637038891
$ws.Range("G3").Value = 11.37596266666667
$ws.Range("H3").Value = 34.127888
$ws.Range("I3").Value = 0.05604480707695051
$ws.Range("J3").Value = 0.05604480707695052
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 1963.019019688112
$ws.Range("R3").Value = 17667.17117719301
$ws.Range("S3").Value = 0.01809049093173171
$ws.Range("T3").Value = 0.01809049093173171
$ws.Range("G4").Value = 11.37596266666667
$ws.Range("H4").Value = 34.127888
$ws.Range("I4").Value = 0.05604480707695051
$ws.Range("J4").Value = 0.05604480707695052
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 846.2317625229883
$ws.Range("R4").Value = 7616.085862706896
$ws.Range("S4").Value = 0.007798573458803137
$ws.Range("T4").Value = 0.007798573458803139
$ws.Range("G5").Value = 11.37596266666667
$ws.Range("H5").Value = 34.127888
$ws.Range("I5").Value = 0.05604480707695051
$ws.Range("J5").Value = 0.05604480707695052
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 664.5224718437315
$ws.Range("R5").Value = 5980.702246593583
$ws.Range("S5").Value = 0.0061240047244835
$ws.Range("T5").Value = 0.006124004724483502
$ws.Range("G6").Value = 11.37596266666667
$ws.Range("H6").Value = 34.127888
$ws.Range("I6").Value = 0.05604480707695051
$ws.Range("J6").Value = 0.05604480707695052
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 793.4303493572694
$ws.Range("R6").Value = 7140.873144215424
$ws.Range("S6").Value = 0.007311974258043063
$ws.Range("T6").Value = 0.007311974258043066
$ws.Range("I7").Value = 0.765548861900355
$ws.Range("J7").Value = 0.7655488619003551
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 24782.31237300697
$ws.Range("R7").Value = 223040.8113570627
$ws.Range("S7").Value = 0.2283850501471227
$ws.Range("T7").Value = 0.2283850501471228
$ws.Range("I8").Value = 0.765548861900355
$ws.Range("J8").Value = 0.7655488619003551
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("S8").Value = 0.2471086165929837
$ws.Range("T8").Value = 0.2471086165929837
$ws.Range("I9").Value = 0.765548861900355
$ws.Range("J9").Value = 0.7655488619003551
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 11559.17553278256
$ws.Range("R9").Value = 104032.5797950431
$ws.Range("S9").Value = 0.1065252848071344
$ws.Range("T9").Value = 0.1065252848071345
$ws.Range("I10").Value = 0.765548861900355
$ws.Range("J10").Value = 0.7655488619003551
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 9077.101850465679
$ws.Range("R10").Value = 81693.9166541911
$ws.Range("S10").Value = 0.08365136917437729
$ws.Range("T10").Value = 0.08365136917437731
$ws.Range("I11").Value = 0.765548861900355
$ws.Range("J11").Value = 0.7655488619003551
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 10837.93008893178
$ws.Range("R11").Value = 97541.37080038601
$ws.Range("S11").Value = 0.09987854117873676
$ws.Range("T11").Value = 0.0998785411787368
$ws.Range("G12").Value = 11.89345866666667
$ws.Range("H12").Value = 35.680376
$ws.Range("I12").Value = 0.05859430238850571
$ws.Range("J12").Value = 0.05859430238850571
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 1896.812048633657
$ws.Range("R12").Value = 17071.30843770292
$ws.Range("S12").Value = 0.01748035083758818
$ws.Range("T12").Value = 0.01748035083758818
$ws.Range("G13").Value = 11.89345866666667
$ws.Range("H13").Value = 35.680376
$ws.Range("I13").Value = 0.05859430238850571
$ws.Range("J13").Value = 0.05859430238850571
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 2052.317351651624
$ws.Range("R13").Value = 18470.85616486462
$ws.Range("S13").Value = 0.01891343286372651
$ws.Range("T13").Value = 0.01891343286372651
$ws.Range("G14").Value = 11.89345866666667
$ws.Range("H14").Value = 35.680376
$ws.Range("I14").Value = 0.05859430238850571
$ws.Range("J14").Value = 0.05859430238850571
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 884.7271026546658
$ws.Range("R14").Value = 7962.543923891993
$ws.Range("S14").Value = 0.008153332936210891
$ws.Range("T14").Value = 0.008153332936210895
$ws.Range("G15").Value = 11.89345866666667
$ws.Range("H15").Value = 35.680376
$ws.Range("I15").Value = 0.05859430238850571
$ws.Range("J15").Value = 0.05859430238850571
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 694.7518011027743
$ws.Range("R15").Value = 6252.766209924968
$ws.Range("S15").Value = 0.006402587561098058
$ws.Range("T15").Value = 0.00640258756109806
$ws.Range("G16").Value = 11.89345866666667
$ws.Range("H16").Value = 35.680376
$ws.Range("I16").Value = 0.05859430238850571
$ws.Range("J16").Value = 0.05859430238850571
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 829.5237371524055
$ws.Range("R16").Value = 7465.713634371649
$ws.Range("S16").Value = 0.007644598189882055
$ws.Range("T16").Value = 0.007644598189882058
$ws.Range("G17").Value = 23.69116533333333
$ws.Range("H17").Value = 71.07349600000001
$ws.Range("I17").Value = 0.1167168730630039
$ws.Range("J17").Value = 0.1167168730630039
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 3778.353220025373
$ws.Range("R17").Value = 34005.17898022835
$ws.Range("S17").Value = 0.03481997065652896
$ws.Range("T17").Value = 0.03481997065652897
$ws.Range("G18").Value = 23.69116533333333
$ws.Range("H18").Value = 71.07349600000001
$ws.Range("I18").Value = 0.1167168730630039
$ws.Range("J18").Value = 0.1167168730630039
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 4088.111882098504
$ws.Range("R18").Value = 36793.00693888654
$ws.Range("S18").Value = 0.0376745972348031
$ws.Range("T18").Value = 0.0376745972348031
$ws.Range("G19").Value = 23.69116533333333
$ws.Range("H19").Value = 71.07349600000001
$ws.Range("I19").Value = 0.1167168730630039
$ws.Range("J19").Value = 0.1167168730630039
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 1762.331433716337
$ws.Range("R19").Value = 15860.98290344703
$ws.Range("S19").Value = 0.01624102492161105
$ws.Range("T19").Value = 0.01624102492161106
$ws.Range("G20").Value = 23.69116533333333
$ws.Range("H20").Value = 71.07349600000001
$ws.Range("I20").Value = 0.1167168730630039
$ws.Range("J20").Value = 0.1167168730630039
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 1383.910286053903
$ws.Range("R20").Value = 12455.19257448513
$ws.Range("S20").Value = 0.01275362909329634
$ws.Range("T20").Value = 0.01275362909329635
$ws.Range("G21").Value = 23.69116533333333
$ws.Range("H21").Value = 71.07349600000001
$ws.Range("I21").Value = 0.1167168730630039
$ws.Range("J21").Value = 0.1167168730630039
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 1652.369134630379
$ws.Range("R21").Value = 14871.32221167341
$ws.Range("S21").Value = 0.01522765115676442
$ws.Range("T21").Value = 0.01522765115676443
$ws.Range("G22").Value = 0.628254
$ws.Range("H22").Value = 1.884762
$ws.Range("I22").Value = 0.003095155571184698
$ws.Range("J22").Value = 0.003095155571184698
$ws.Range("M22").Value = 159.4836373333333
$ws.Range("N22").Value = 478.450912
$ws.Range("O22").Value = 0.2983285084902258
$ws.Range("P22").Value = 0.2983285084902258
$ws.Range("Q22").Value = 100.196233089216
$ws.Range("R22").Value = 901.7660978029439
$ws.Range("S22").Value = 0.0009233731450967438
$ws.Range("T22").Value = 0.0009233731450967439
$ws.Range("G23").Value = 0.628254
$ws.Range("H23").Value = 1.884762
$ws.Range("I23").Value = 0.003095155571184698
$ws.Range("J23").Value = 0.003095155571184698
$ws.Range("O23").Value = 0.3227862111630279
$ws.Range("P23").Value = 0.3227862111630279
$ws.Range("Q23").Value = 108.410565974238
$ws.Range("R23").Value = 975.6950937681419
$ws.Range("S23").Value = 0.0009990735397828461
$ws.Range("T23").Value = 0.0009990735397828463
$ws.Range("G24").Value = 0.628254
$ws.Range("H24").Value = 1.884762
$ws.Range("I24").Value = 0.003095155571184698
$ws.Range("J24").Value = 0.003095155571184698
$ws.Range("M24").Value = 74.38770566666666
$ws.Range("N24").Value = 223.163117
$ws.Range("O24").Value = 0.1391489036280481
$ws.Range("P24").Value = 0.1391489036280482
$ws.Range("Q24").Value = 46.734373635906
$ws.Range("R24").Value = 420.6093627231539
$ws.Range("S24").Value = 0.0004306875042885958
$ws.Range("T24").Value = 0.0004306875042885959
$ws.Range("G25").Value = 0.628254
$ws.Range("H25").Value = 1.884762
$ws.Range("I25").Value = 0.003095155571184698
$ws.Range("J25").Value = 0.003095155571184698
$ws.Range("M25").Value = 58.41461433333333
$ws.Range("N25").Value = 175.243843
$ws.Range("O25").Value = 0.1092697975759847
$ws.Range("P25").Value = 0.1092697975759848
$ws.Range("Q25").Value = 36.699215113374
$ws.Range("R25").Value = 330.292936020366
$ws.Range("S25").Value = 0.0003382070227295334
$ws.Range("T25").Value = 0.0003382070227295335
$ws.Range("G26").Value = 0.628254
$ws.Range("H26").Value = 1.884762
$ws.Range("I26").Value = 0.003095155571184698
$ws.Range("J26").Value = 0.003095155571184698
$ws.Range("M26").Value = 69.746216
$ws.Range("N26").Value = 209.238648
$ws.Range("O26").Value = 0.1304665791427133
$ws.Range("P26").Value = 0.1304665791427133
$ws.Range("Q26").Value = 43.818339186864
$ws.Range("R26").Value = 394.365052681776
$ws.Range("S26").Value = 0.0004038143592869784
$ws.Range("T26").Value = 0.0004038143592869785
